$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TrancheDetails")

$ws.Range("A5").Value = "Frac Total Gas"
$ws.Range("B5").Value = 0.00029534195905716725
$ws.Range("C5").Value = 0.0036453349502072671
$ws.Range("D5").Value = 0.0066109069432481451
$ws.Range("E5").Value = 0.015144015324681041
$ws.Range("F5").Value = 0.056319634110065388
$ws.Range("G5").Value = 0.078406143359654346
$ws.Range("H5").Value = 0.25340652577454159
$ws.Range("I5").Value = 0.1311533569343277
$ws.Range("J5").Value = 0.4145126613663031
$ws.Range("K5").Value = 0.040506079277914275

$ws.Range("A6").Select()
